# Daily attendance processing - 2025-12-03 11:26:14
#
# Normalizes the "Recorded By" (column G) text so that "System" is listed
# first among the recorder names, for the specific value patterns that were
# recorded with "System" trailing at the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

$usedRange = $ws.UsedRange
$maxRow = $usedRange.Rows.Count + $usedRange.Row - 1
if ($maxRow -gt $lastRow) { $lastRow = $maxRow }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($null -eq $val) { continue }

    switch ($val) {
        "backup@backdoor.com, system, System" {
            $cell.Value = "backup@backdoor.com, System, system"
        }
        "dnasr281@gmail.com, System" {
            $cell.Value = "System, dnasr281@gmail.com"
        }
        "admin@admin.com, System" {
            $cell.Value = "System, admin@admin.com"
        }
    }
}
